$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{row=2; B=1; C=0.01180555555555556; D=0.03611111111111111; E="F36"; F=1},
    @{row=3; B=2; C=0.02638888888888889; D=0.04722222222222222; E="E22"; F=2},
    @{row=4; B=3; C=0.03958333333333333; D=0.06666666666666667; E="D33"; F=2},
    @{row=5; B=4; C=0.02430555555555556; D=0.04861111111111111; E="D45"; F=2},
    @{row=6; B=5; C=0.02083333333333333; D=0.04722222222222222; E="F33"; F=3},
    @{row=7; B=6; C=0.03472222222222222; D=0.05694444444444444; E="A6"; F=2},
    @{row=8; B=7; C=0.00625; D=0.03125; E="A12"; F=2},
    @{row=9; B=8; C=0.02152777777777778; D=0.04305555555555556; E="F40"; F=3},
    @{row=10; B=9; C=0.03402777777777777; D=0.05486111111111111; E="F60"; F=2},
    @{row=11; B=10; C=0.0125; D=0.03611111111111111; E="F54"; F=3},
    @{row=12; B=11; C=0.03263888888888889; D=0.05416666666666667; E="F51"; F=1},
    @{row=13; B=12; C=0.02916666666666667; D=0.05208333333333334; E="A10"; F=2},
    @{row=14; B=13; C=0.01875; D=0.04444444444444445; E="C1"; F=1},
    @{row=15; B=14; C=0.01944444444444444; D=0.04166666666666666; E="A18"; F=1},
    @{row=16; B=15; C=0.01527777777777778; D=0.03819444444444445; E="C16"; F=3},
    @{row=17; B=16; C=0.04861111111111111; D=0.07222222222222222; E="E1"; F=1},
    @{row=18; B=17; C=0.05138888888888889; D=0.07222222222222222; E="E28"; F=2},
    @{row=19; B=18; C=0.06319444444444444; D=0.08402777777777778; E="D32"; F=3},
    @{row=20; B=19; C=0.05; D=0.07152777777777777; E="E7"; F=1},
    @{row=21; B=20; C=0.04236111111111111; D=0.06875000000000001; E="A1"; F=1}
)

foreach ($item in $data) {
    $r = $item.row
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
}

Write-Host "Applied schedule_8 instance differences"
